# Update the "want to go" counts (column F) for a handful of events.
# Sheet 1 = 展览 (Exhibition), Sheet 4 = 全部类型 (All types) — the same
# events are listed on both sheets, so every value is bumped in both places.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F3").Value = 601
$wsExhibition.Range("F7").Value = 14977
$wsExhibition.Range("F13").Value = 8749
$wsExhibition.Range("F31").Value = 31
$wsExhibition.Range("F36").Value = 273
$wsExhibition.Range("F38").Value = 110

# Sheet 4: 全部类型 (mirrors the same rows, offset by one due to extra row)
$wsAllTypes = $wb.Worksheets.Item(4)
$wsAllTypes.Range("F3").Value = 601
$wsAllTypes.Range("F7").Value = 14977
$wsAllTypes.Range("F13").Value = 8749
$wsAllTypes.Range("F32").Value = 31
$wsAllTypes.Range("F39").Value = 273
$wsAllTypes.Range("F41").Value = 110
